$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.493.47"
$ws.Range("E2").Value = "  +5.04%  "
$ws.Range("D3").Value = "3.323.49"
$ws.Range("E3").Value = "  +4.60%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'552.44"
$ws.Range("E5").Value = "  +3.28%  "
$ws.Range("D6").Value = "'151.50"
$ws.Range("E6").Value = "  +5.10%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'0.526"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "'7.50"
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("D10").Value = "'0.118"
$ws.Range("E10").Value = "  +4.50%  "
$ws.Range("D11").Value = "'0.436"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "3.900.16"
$ws.Range("E12").Value = "  +4.70%  "
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000181"
$ws.Range("E14").Value = "  +4.93%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'27.02"
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("D16").Value = "62.435.54"
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("D17").Value = "3.327.88"
$ws.Range("E17").Value = "  +4.46%  "
$ws.Range("D18").Value = "'6.47"
$ws.Range("E18").Value = "  +4.48%  "
$ws.Range("D19").Value = "'13.74"
$ws.Range("E19").Value = "  +5.90%  "
$ws.Range("D20").Value = "'8.49"
$ws.Range("E20").Value = "  +4.02%  "
$ws.Range("D21").Value = "'385.06"
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'0.536"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").Value = "'70.76"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "'0.177"
$ws.Range("E25").Value = "  +4.30%  "
$ws.Range("D26").Value = "'8.88"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").Value = "0.0₃0973"
$ws.Range("E27").Value = "  +8.75%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  +3.61%  "
$ws.Range("D30").Value = "'6.36"
$ws.Range("E30").Value = "  +3.90%  "
$ws.Range("E31").Value = "  +11.55%  "
$ws.Range("D32").Value = "'22.93"
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("D33").Value = "'5.54"
$ws.Range("E33").Value = "  +3.50%  "
$ws.Range("D34").Value = "'6.73"
$ws.Range("E34").Value = "  +4.83%  "
$ws.Range("E35").Value = "  +11.15%  "
$ws.Range("D36").Value = "'159.47"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("E37").Value = "  +12.45%  "
$ws.Range("D38").Value = "'27.06"
$ws.Range("E38").Value = "  +5.86%  "
$ws.Range("D39").Value = "2.852.70"
$ws.Range("E39").Value = "  +4.42%  "
$ws.Range("D40").Value = "'0.0731"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("D41").Value = "'0.0315"
$ws.Range("E41").Value = "  +8.54%  "
$ws.Range("D42").Value = "'4.32"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("E43").Value = "  +3.42%  "
$ws.Range("D44").Value = "'40.71"
$ws.Range("E44").Value = "  +3.65%  "
$ws.Range("D45").Value = "'1.04"
$ws.Range("E45").Value = "  +4.27%  "
$ws.Range("D46").Value = "'21.97"
$ws.Range("E46").Value = "  +7.40%  "
$ws.Range("D47").Value = "3.370.59"
$ws.Range("E47").Value = "  +4.64%  "
$ws.Range("E48").Value = "  +4.27%  "
$ws.Range("D49").Value = "'6.28"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").Value = "'0.807"
$ws.Range("E50").Value = "  +5.05%  "
$ws.Range("D51").Value = "'282.81"
$ws.Range("E51").Value = "  +8.23%  "
